# New weekly price record for Jengibre (Terminal La Palmera de La Serena)
# is inserted as row 44, pushing the existing rows 44-106 down to 45-107.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44 (shifts rows 44:106 down to 45:107,
# carrying their formatting/number-format along, same as Excel's
# EntireRow.Insert / Rows.Insert behaviour).
$ws.Rows.Item(44).Insert()

# Fill the newly inserted row 44 with the new record's data.
$ws.Cells.Item(44, 1).Value = 8
$ws.Cells.Item(44, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(44, 3).Value = "Coquimbo"
$ws.Cells.Item(44, 4).Value = 44967
$ws.Cells.Item(44, 5).Value = 4
$ws.Cells.Item(44, 6).Value = 100114007
$ws.Cells.Item(44, 7).Value = "Jengibre"
$ws.Cells.Item(44, 8).Value = "Sin especificar"
$ws.Cells.Item(44, 9).Value = "Primera"
$ws.Cells.Item(44, 10).Value = 400
$ws.Cells.Item(44, 11).Value = 22000
$ws.Cells.Item(44, 12).Value = 23000
$ws.Cells.Item(44, 13).Value = 22500
$ws.Cells.Item(44, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(44, 15).Value = "Perú"
$ws.Cells.Item(44, 16).Value = 1731
$ws.Cells.Item(44, 17).Value = 13
$ws.Cells.Item(44, 18).Value = "Hortaliza"
